# Update the "Presented By" textbox on slide 1 (shape "TextBox 3") with the
# student / college details, splitting the text into the same run layout as
# the authored edit (so that distinct formatting -- e.g. the superscript
# "rd" -- lands on its own run).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$enDash = [char]0x2013

# Paragraph text pieces (concatenated with CR to form paragraph breaks).
$p1r1 = "Presented By"
$p1r2 = ":"

$p2r1 = "                 NISHA A"

$p3r1 = " "
$p3r2 = "               3"
$p3r3 = "rd"
$p3r4 = " " + $enDash
$p3r5 = "B.Tech"
$p3r6 = "-IT"

$p4r1 = " "
$p4r2 = "               "
$p4r3 = "Mookambigai"
$p4r4 = " college "
$p4r5 = "ofengineering"
$p4r6 = " "

$para1 = $p1r1 + $p1r2
$para2 = $p2r1
$para3 = $p3r1 + $p3r2 + $p3r3 + $p3r4 + $p3r5 + $p3r6
$para4 = $p4r1 + $p4r2 + $p4r3 + $p4r4 + $p4r5 + $p4r6

$tr.Text = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4

# Common look for every run in this textbox: 20pt bold, Arial, accent1 75% luma.
function Set-RunLook($range) {
    $range.Font.Name = "Arial"
    $range.Font.Size = 20
    $range.Font.Bold = $true
}

$pos = 1

# Paragraph 1: "Presented By" + ":"
$r = $tr.Characters($pos, $p1r1.Length); Set-RunLook $r
$pos += $p1r1.Length
$r = $tr.Characters($pos, $p1r2.Length); Set-RunLook $r
$pos += $p1r2.Length
$pos += 1 # CR

# Paragraph 2: "                 NISHA A"
$r = $tr.Characters($pos, $p2r1.Length); Set-RunLook $r
$pos += $p2r1.Length
$pos += 1 # CR

# Paragraph 3: " " + "               3" + "rd" (superscript) + " –" + "B.Tech" + "-IT"
$r = $tr.Characters($pos, $p3r1.Length); Set-RunLook $r
$pos += $p3r1.Length
$r = $tr.Characters($pos, $p3r2.Length); Set-RunLook $r
$pos += $p3r2.Length
$r = $tr.Characters($pos, $p3r3.Length); Set-RunLook $r; $r.Font.BaselineOffset = 0.3
$pos += $p3r3.Length
$r = $tr.Characters($pos, $p3r4.Length); Set-RunLook $r
$pos += $p3r4.Length
$r = $tr.Characters($pos, $p3r5.Length); Set-RunLook $r
$pos += $p3r5.Length
$r = $tr.Characters($pos, $p3r6.Length); Set-RunLook $r
$pos += $p3r6.Length
$pos += 1 # CR

# Paragraph 4: " " + spaces + "Mookambigai" + " college " + "ofengineering" + " "
$r = $tr.Characters($pos, $p4r1.Length); Set-RunLook $r
$pos += $p4r1.Length
$r = $tr.Characters($pos, $p4r2.Length); Set-RunLook $r
$pos += $p4r2.Length
$r = $tr.Characters($pos, $p4r3.Length); Set-RunLook $r
$pos += $p4r3.Length
$r = $tr.Characters($pos, $p4r4.Length); Set-RunLook $r
$pos += $p4r4.Length
$r = $tr.Characters($pos, $p4r5.Length); Set-RunLook $r
$pos += $p4r5.Length
$r = $tr.Characters($pos, $p4r6.Length); Set-RunLook $r
$pos += $p4r6.Length

# Re-apply the accent1/lumMod75 text color across the whole box (color set
# last so every run -- including the ones just split out above -- carries it).
$tr.Font.Color.ObjectThemeColor = 5   # msoThemeColorAccent1
$tr.Font.Color.Brightness = -0.25     # lumMod 75% / lumOff 0%
